$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36
$ws.Range("B36").Value2 = 6864629.0
$ws.Range("F36").Value2 = 'Borac Banja Luka'
$ws.Range("G36").Value2 = 'NK Posusje'
$ws.Range("H36").Value2 = 1.0
$ws.Range("I36").Value2 = 0.0
$ws.Range("J36").Value2 = 'H'
$ws.Range("K36").Value2 = 1.363
$ws.Range("L36").Value2 = 4.5
$ws.Range("M36").Value2 = 6.5
$ws.Range("N36").Value2 = 1.363
$ws.Range("O36").Value2 = 4.2
$ws.Range("P36").Value2 = 6.5
$ws.Range("Q36").Value2 = -1.25
$ws.Range("R36").Value2 = 1.95
$ws.Range("S36").Value2 = 1.85
$ws.Range("T36").Value2 = 2.5
$ws.Range("U36").Value2 = 1.925
$ws.Range("V36").Value2 = 1.875
$ws.Range("W36").Value2 = 0.363
$ws.Range("X36").Value2 = -1.0
$ws.Range("Y36").Value2 = -1.0
$ws.Range("Z36").Value2 = -0.5
$ws.Range("AA36").Value2 = 0.425
$ws.Range("AB36").Value2 = -1.0
$ws.Range("AC36").Value2 = 0.875

# Row 37
$ws.Range("B37").Value2 = 6865299.0
$ws.Range("F37").Value2 = 'Siroki Brijeg'
$ws.Range("G37").Value2 = 'Zvijezda 09'
$ws.Range("H37").Value2 = 2.0
$ws.Range("I37").Value2 = 1.0
$ws.Range("J37").Value2 = 'H'
$ws.Range("K37").Value2 = 1.25
$ws.Range("L37").Value2 = 5.5
$ws.Range("M37").Value2 = 8.0
$ws.Range("N37").Value2 = 1.4
$ws.Range("O37").Value2 = 4.75
$ws.Range("P37").Value2 = 5.75
$ws.Range("Q37").Value2 = -1.25
$ws.Range("R37").Value2 = 1.9
$ws.Range("S37").Value2 = 1.9
$ws.Range("T37").Value2 = 2.75
$ws.Range("U37").Value2 = 1.85
$ws.Range("V37").Value2 = 1.95
$ws.Range("W37").Value2 = 0.3999999999999999
$ws.Range("X37").Value2 = -1.0
$ws.Range("Y37").Value2 = -1.0
$ws.Range("Z37").Value2 = -0.5
$ws.Range("AA37").Value2 = 0.45
$ws.Range("AB37").Value2 = 0.425
$ws.Range("AC37").Value2 = -0.5

# Row 87
$ws.Range("B87").Value2 = 7505495.0
$ws.Range("F87").Value2 = 'Sloga'
$ws.Range("G87").Value2 = 'Zvijezda 09'
$ws.Range("H87").Value2 = 1.0
$ws.Range("I87").Value2 = 0.0
$ws.Range("J87").Value2 = 'H'
$ws.Range("K87").Value2 = 1.444
$ws.Range("L87").Value2 = 4.2
$ws.Range("M87").Value2 = 5.5
$ws.Range("N87").Value2 = 1.5
$ws.Range("O87").Value2 = 4.2
$ws.Range("P87").Value2 = 5.25
$ws.Range("Q87").Value2 = -1.0
$ws.Range("R87").Value2 = 1.8
$ws.Range("S87").Value2 = 2.0
$ws.Range("T87").Value2 = 2.75
$ws.Range("U87").Value2 = 1.775
$ws.Range("V87").Value2 = 2.025
$ws.Range("W87").Value2 = 0.5
$ws.Range("X87").Value2 = -1.0
$ws.Range("Y87").Value2 = -1.0
$ws.Range("Z87").Value2 = 0.0
$ws.Range("AA87").Value2 = 0.0
$ws.Range("AB87").Value2 = -1.0
$ws.Range("AC87").Value2 = 1.025

# Row 88
$ws.Range("B88").Value2 = 7505497.0
$ws.Range("F88").Value2 = 'Zeljeznicar'
$ws.Range("G88").Value2 = 'NK Posusje'
$ws.Range("H88").Value2 = 1.0
$ws.Range("I88").Value2 = 1.0
$ws.Range("J88").Value2 = 'D'
$ws.Range("K88").Value2 = 1.65
$ws.Range("L88").Value2 = 3.4
$ws.Range("M88").Value2 = 4.75
$ws.Range("N88").Value2 = 1.8
$ws.Range("O88").Value2 = 3.2
$ws.Range("P88").Value2 = 4.2
$ws.Range("Q88").Value2 = -0.5
$ws.Range("R88").Value2 = 1.825
$ws.Range("S88").Value2 = 1.975
$ws.Range("T88").Value2 = 2.0
$ws.Range("U88").Value2 = 1.75
$ws.Range("V88").Value2 = 2.05
$ws.Range("W88").Value2 = -1.0
$ws.Range("X88").Value2 = 2.2
$ws.Range("Y88").Value2 = -1.0
$ws.Range("Z88").Value2 = -1.0
$ws.Range("AA88").Value2 = 0.9750000000000001
$ws.Range("AB88").Value2 = 0.0
$ws.Range("AC88").Value2 = 0.0

# Row 99
$ws.Range("B99").Value2 = 6865343.0
$ws.Range("F99").Value2 = 'Sloga'
$ws.Range("G99").Value2 = 'NK Posusje'
$ws.Range("H99").Value2 = 1.0
$ws.Range("I99").Value2 = 0.0
$ws.Range("J99").Value2 = 'H'
$ws.Range("K99").Value2 = 1.909
$ws.Range("L99").Value2 = 3.3
$ws.Range("M99").Value2 = 3.5
$ws.Range("N99").Value2 = 2.2
$ws.Range("O99").Value2 = 2.8
$ws.Range("P99").Value2 = 3.3
$ws.Range("Q99").Value2 = -0.25
$ws.Range("R99").Value2 = 1.95
$ws.Range("S99").Value2 = 1.85
$ws.Range("T99").Value2 = 1.75
$ws.Range("U99").Value2 = 1.875
$ws.Range("V99").Value2 = 1.925
$ws.Range("W99").Value2 = 1.2
$ws.Range("X99").Value2 = -1.0
$ws.Range("Y99").Value2 = -1.0
$ws.Range("Z99").Value2 = 0.95
$ws.Range("AA99").Value2 = -1.0
$ws.Range("AB99").Value2 = -1.0
$ws.Range("AC99").Value2 = 0.925

# Row 100
$ws.Range("B100").Value2 = 6864639.0
$ws.Range("F100").Value2 = 'Zvijezda 09'
$ws.Range("G100").Value2 = 'Borac Banja Luka'
$ws.Range("H100").Value2 = 1.0
$ws.Range("I100").Value2 = 2.0
$ws.Range("J100").Value2 = 'A'
$ws.Range("K100").Value2 = 11.0
$ws.Range("L100").Value2 = 6.0
$ws.Range("M100").Value2 = 1.2
$ws.Range("N100").Value2 = 10.0
$ws.Range("O100").Value2 = 6.5
$ws.Range("P100").Value2 = 1.181
$ws.Range("Q100").Value2 = 2.0
$ws.Range("R100").Value2 = 1.825
$ws.Range("S100").Value2 = 1.975
$ws.Range("T100").Value2 = 3.0
$ws.Range("U100").Value2 = 1.9
$ws.Range("V100").Value2 = 1.9
$ws.Range("W100").Value2 = -1.0
$ws.Range("X100").Value2 = -1.0
$ws.Range("Y100").Value2 = 0.181
$ws.Range("Z100").Value2 = 0.825
$ws.Range("AA100").Value2 = -1.0
$ws.Range("AB100").Value2 = 0.0
$ws.Range("AC100").Value2 = 0.0

# Row 111
$ws.Range("B111").Value2 = 6865354.0
$ws.Range("F111").Value2 = 'NK Igman Konjic'
$ws.Range("G111").Value2 = 'GOSK Gabela'
$ws.Range("H111").Value2 = 1.0
$ws.Range("I111").Value2 = 2.0
$ws.Range("J111").Value2 = 'A'
$ws.Range("K111").Value2 = 1.8
$ws.Range("L111").Value2 = 3.25
$ws.Range("M111").Value2 = 4.0
$ws.Range("N111").Value2 = 2.25
$ws.Range("O111").Value2 = 3.1
$ws.Range("P111").Value2 = 2.9
$ws.Range("Q111").Value2 = -0.25
$ws.Range("R111").Value2 = 1.975
$ws.Range("S111").Value2 = 1.825
$ws.Range("T111").Value2 = 2.25
$ws.Range("U111").Value2 = 1.875
$ws.Range("V111").Value2 = 1.925
$ws.Range("W111").Value2 = -1.0
$ws.Range("X111").Value2 = -1.0
$ws.Range("Y111").Value2 = 1.9
$ws.Range("Z111").Value2 = -1.0
$ws.Range("AA111").Value2 = 0.825
$ws.Range("AB111").Value2 = 0.875
$ws.Range("AC111").Value2 = -1.0

# Row 112
$ws.Range("B112").Value2 = 6865352.0
$ws.Range("F112").Value2 = 'NK Posusje'
$ws.Range("G112").Value2 = 'Zvijezda 09'
$ws.Range("H112").Value2 = 2.0
$ws.Range("I112").Value2 = 0.0
$ws.Range("J112").Value2 = 'H'
$ws.Range("K112").Value2 = 1.615
$ws.Range("L112").Value2 = 3.5
$ws.Range("M112").Value2 = 4.75
$ws.Range("N112").Value2 = 1.5
$ws.Range("O112").Value2 = 3.6
$ws.Range("P112").Value2 = 5.75
$ws.Range("Q112").Value2 = -1.0
$ws.Range("R112").Value2 = 1.9
$ws.Range("S112").Value2 = 1.9
$ws.Range("T112").Value2 = 2.25
$ws.Range("U112").Value2 = 1.85
$ws.Range("V112").Value2 = 1.95
$ws.Range("W112").Value2 = 0.5
$ws.Range("X112").Value2 = -1.0
$ws.Range("Y112").Value2 = -1.0
$ws.Range("Z112").Value2 = 0.8999999999999999
$ws.Range("AA112").Value2 = -1.0
$ws.Range("AB112").Value2 = -0.5
$ws.Range("AC112").Value2 = 0.475

# Row 125
$ws.Range("B125").Value2 = 6865362.0
$ws.Range("F125").Value2 = 'NK Posusje'
$ws.Range("G125").Value2 = 'Velez Mostar'
$ws.Range("H125").Value2 = 1.0
$ws.Range("I125").Value2 = 0.0
$ws.Range("J125").Value2 = 'H'
$ws.Range("K125").Value2 = 3.4
$ws.Range("L125").Value2 = 2.9
$ws.Range("M125").Value2 = 2.15
$ws.Range("N125").Value2 = 2.625
$ws.Range("O125").Value2 = 2.8
$ws.Range("P125").Value2 = 2.625
$ws.Range("Q125").Value2 = 0.0
$ws.Range("R125").Value2 = 1.9
$ws.Range("S125").Value2 = 1.9
$ws.Range("T125").Value2 = 1.75
$ws.Range("U125").Value2 = 1.8
$ws.Range("V125").Value2 = 2.0
$ws.Range("W125").Value2 = 1.625
$ws.Range("X125").Value2 = -1.0
$ws.Range("Y125").Value2 = -1.0
$ws.Range("Z125").Value2 = 0.8999999999999999
$ws.Range("AA125").Value2 = -1.0
$ws.Range("AB125").Value2 = -1.0
$ws.Range("AC125").Value2 = 1.0

# Row 126
$ws.Range("B126").Value2 = 6865364.0
$ws.Range("F126").Value2 = 'Zeljeznicar'
$ws.Range("G126").Value2 = 'FK Sarajevo'
$ws.Range("H126").Value2 = 3.0
$ws.Range("I126").Value2 = 0.0
$ws.Range("J126").Value2 = 'H'
$ws.Range("K126").Value2 = 3.25
$ws.Range("L126").Value2 = 3.0
$ws.Range("M126").Value2 = 2.2
$ws.Range("N126").Value2 = 2.7
$ws.Range("O126").Value2 = 2.8
$ws.Range("P126").Value2 = 2.7
$ws.Range("Q126").Value2 = 0.0
$ws.Range("R126").Value2 = 1.875
$ws.Range("S126").Value2 = 1.925
$ws.Range("T126").Value2 = 2.0
$ws.Range("U126").Value2 = 2.05
$ws.Range("V126").Value2 = 1.75
$ws.Range("W126").Value2 = 1.7
$ws.Range("X126").Value2 = -1.0
$ws.Range("Y126").Value2 = -1.0
$ws.Range("Z126").Value2 = 0.875
$ws.Range("AA126").Value2 = -1.0
$ws.Range("AB126").Value2 = 1.05
$ws.Range("AC126").Value2 = -1.0

# New row 127
$ws.Range("A126").Copy() | Out-Null
$ws.Range("A127").PasteSpecial(-4122) | Out-Null
$ws.Range("E126").Copy() | Out-Null
$ws.Range("E127").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A127").Value2 = 125.0
$ws.Range("B127").Value2 = 6864643.0
$ws.Range("C127").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D127").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E127").Value2 = 45355.57291666666
$ws.Range("F127").Value2 = 'Borac Banja Luka'
$ws.Range("G127").Value2 = 'GOSK Gabela'
$ws.Range("H127").Value2 = 3.0
$ws.Range("I127").Value2 = 0.0
$ws.Range("J127").Value2 = 'H'
$ws.Range("K127").Value2 = 1.181
$ws.Range("L127").Value2 = 6.5
$ws.Range("M127").Value2 = 11.0
$ws.Range("N127").Value2 = 1.222
$ws.Range("O127").Value2 = 6.0
$ws.Range("P127").Value2 = 9.0
$ws.Range("Q127").Value2 = -1.75
$ws.Range("R127").Value2 = 1.85
$ws.Range("S127").Value2 = 1.95
$ws.Range("T127").Value2 = 3.0
$ws.Range("U127").Value2 = 1.95
$ws.Range("V127").Value2 = 1.85
$ws.Range("W127").Value2 = 0.222
$ws.Range("X127").Value2 = -1.0
$ws.Range("Y127").Value2 = -1.0
$ws.Range("Z127").Value2 = 0.8500000000000001
$ws.Range("AA127").Value2 = -1.0
$ws.Range("AB127").Value2 = 0.0
$ws.Range("AC127").Value2 = 0.0
